$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 77 (pushes the existing row 77 and everything
# below it down by one, i.e. old row 77 becomes row 78, old row 78 becomes
# row 79, ..., old row 115 becomes row 116). The new row 77 holds a new
# pineapple price record (same as the old row 77 data but with an updated
# date of 44518 instead of 44446).
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with a duplicate of the (now shifted
# down) old row 77 data, but with the updated date.
$ws.Range("A77").Value = 11
$ws.Range("B77").Value = "Vega Monumental Concepción"
$ws.Range("C77").Value = "Bíobío"
$ws.Range("D77").Value = 44518
$ws.Range("D77").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108005
$ws.Range("J77").Value = "Piña"
$ws.Range("K77").Value = "Caramelo"
$ws.Range("L77").Value = "Segunda"
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 18000
$ws.Range("O77").Value = 19000
$ws.Range("P77").Value = 18500
$ws.Range("Q77").Value = "$/caja 14 unidades"
$ws.Range("R77").Value = "Ecuador"
$ws.Range("S77").Value = 1321
$ws.Range("T77").Value = 14
